# Completed implementation of RQ3 for both regular and micro-clones.
#
# "RQ3" (regular clones) lives on the worksheet tab named "RQ2" and
# "RQ3 micro" (micro-clones) lives on the worksheet tab named "RQ3" -
# both previously had their first data row (row 3 / row 4 respectively)
# empty, which made the %-formulas in column D/G evaluate to #DIV/0!.
# Fill in the missing raw counts so the percentages compute correctly;
# the dependent SUM row (row 9) recalculates automatically.

$wb = $excel.ActiveWorkbook

# --- Worksheet "RQ2": first project's regular-clone numbers (row 3) ---
$wsRQ2 = $wb.Worksheets.Item("RQ2")
$wsRQ2.Range("B3").Value = 83
$wsRQ2.Range("C3").Value = 86
$wsRQ2.Range("D3").Formula = "=B3/C3*100"
$wsRQ2.Range("E3").Value = 1225
$wsRQ2.Range("F3").Value = 13401
$wsRQ2.Range("G3").Formula = "=E3/F3*100"

# --- Worksheet "RQ3": first project's micro-clone numbers (row 4) ---
$wsRQ3 = $wb.Worksheets.Item("RQ3")
$wsRQ3.Range("B4").Value = 12
$wsRQ3.Range("C4").Value = 12
$wsRQ3.Range("D4").Formula = "=B4/C4*100"
$wsRQ3.Range("E4").Value = 72
$wsRQ3.Range("F4").Value = 72
$wsRQ3.Range("G4").Formula = "=E4/F4*100"

# --- Leave the UI state matching the author's final save ---
$wsRQ2.Range("F3").Select()
$wsRQ3.Activate()
$wsRQ3.Range("B4").Select()
